# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties
# Header row (row 1) gets the same style as the other header cells.
# Data rows (2-46) get the team's W/L/T record: 86-75-0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the bold/border/center style
# used by the rest of row 1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-46) ---------------------------------------------------
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 75   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
